$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Razon social (column E) comma-to-period fixes ---
$ws.Range("E93").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E103").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E238").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E278").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E189").Value = "DODERA. JORGE ABELARDO"
$ws.Range("E196").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E209").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E197").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"

# --- Importe (column H) thousands/decimal separator fixes ---
$ws.Range("H2").Value = "'38025.00"
$ws.Range("H3").Value = "'5550.00"
$ws.Range("H4").Value = "'16950.00"
$ws.Range("H5").Value = "'49600.00"
$ws.Range("H6").Value = "'138000.00"
$ws.Range("H7").Value = "'4373.50"
$ws.Range("H8").Value = "'11730.00"
$ws.Range("H9").Value = "'1480000.00"
$ws.Range("H10").Value = "'35000.00"
$ws.Range("H11").Value = "'86000.00"
$ws.Range("H12").Value = "'3089.80"
$ws.Range("H13").Value = "'4666.00"
$ws.Range("H14").Value = "'8100.00"
$ws.Range("H15").Value = "'64000.00"
$ws.Range("H16").Value = "'210000.00"
$ws.Range("H17").Value = "'103666.00"
$ws.Range("H18").Value = "'63810.00"
$ws.Range("H19").Value = "'16653.60"
$ws.Range("H20").Value = "'15544.00"
$ws.Range("H21").Value = "'50680.50"
$ws.Range("H22").Value = "'14564.00"
$ws.Range("H23").Value = "'80555.74"
$ws.Range("H24").Value = "'11378.06"
$ws.Range("H25").Value = "'5349.52"
$ws.Range("H26").Value = "'29930.12"
$ws.Range("H27").Value = "'5635.00"
$ws.Range("H28").Value = "'87144.00"
$ws.Range("H29").Value = "'36163.59"
$ws.Range("H30").Value = "'4900.00"
$ws.Range("H31").Value = "'2800.00"
$ws.Range("H32").Value = "'6433.00"
$ws.Range("H33").Value = "'162050.00"
$ws.Range("H34").Value = "'661.02"
$ws.Range("H35").Value = "'8718.25"
$ws.Range("H36").Value = "'8200.00"
$ws.Range("H37").Value = "'14820.00"
$ws.Range("H38").Value = "'102402.58"
$ws.Range("H39").Value = "'400.57"
$ws.Range("H40").Value = "'150.00"
$ws.Range("H41").Value = "'1127.39"
$ws.Range("H42").Value = "'87624.50"
$ws.Range("H43").Value = "'3087000.00"
$ws.Range("H44").Value = "'6459.51"
$ws.Range("H45").Value = "'4658.00"
$ws.Range("H46").Value = "'2000.00"
$ws.Range("H47").Value = "'2394.40"
$ws.Range("H48").Value = "'597493.02"
$ws.Range("H49").Value = "'13460.62"
$ws.Range("H50").Value = "'96.00"
$ws.Range("H51").Value = "'3801.32"
$ws.Range("H52").Value = "'940.00"
$ws.Range("H53").Value = "'12400.00"
$ws.Range("H54").Value = "'8527.29"
$ws.Range("H55").Value = "'5872.74"
$ws.Range("H56").Value = "'842.80"
$ws.Range("H57").Value = "'8300.00"
$ws.Range("H58").Value = "'7674.96"
$ws.Range("H59").Value = "'11200.00"
$ws.Range("H60").Value = "'301700.74"
$ws.Range("H61").Value = "'3055.21"
$ws.Range("H62").Value = "'17441.27"
$ws.Range("H63").Value = "'5019.99"
$ws.Range("H64").Value = "'3560.00"
$ws.Range("H65").Value = "'2800.00"
$ws.Range("H66").Value = "'9510.43"
$ws.Range("H67").Value = "'10230.00"
$ws.Range("H68").Value = "'27900.00"
$ws.Range("H69").Value = "'26200.00"
$ws.Range("H70").Value = "'6402.66"
$ws.Range("H71").Value = "'4081.08"
$ws.Range("H72").Value = "'300.00"
$ws.Range("H73").Value = "'2492.80"
$ws.Range("H74").Value = "'10480.00"
$ws.Range("H75").Value = "'6009.53"
$ws.Range("H76").Value = "'3890.75"
$ws.Range("H77").Value = "'19559.87"
$ws.Range("H78").Value = "'2950.00"
$ws.Range("H79").Value = "'968.00"
$ws.Range("H80").Value = "'10276.19"
$ws.Range("H81").Value = "'1424.00"
$ws.Range("H82").Value = "'4000.00"
$ws.Range("H83").Value = "'5000.00"
$ws.Range("H84").Value = "'534.00"
$ws.Range("H85").Value = "'5630.00"
$ws.Range("H86").Value = "'17155.72"
$ws.Range("H87").Value = "'5000.00"
$ws.Range("H88").Value = "'30560.00"
$ws.Range("H89").Value = "'2000.00"
$ws.Range("H90").Value = "'1280.00"
$ws.Range("H91").Value = "'756.00"
$ws.Range("H92").Value = "'13254.00"
$ws.Range("H93").Value = "'5021.00"
$ws.Range("H94").Value = "'7548.00"
$ws.Range("H95").Value = "'140.25"
$ws.Range("H96").Value = "'880.00"
$ws.Range("H97").Value = "'1313.00"
$ws.Range("H98").Value = "'406.00"
$ws.Range("H99").Value = "'192.00"
$ws.Range("H100").Value = "'3875.00"
$ws.Range("H101").Value = "'590620.31"
$ws.Range("H102").Value = "'1700.00"
$ws.Range("H103").Value = "'2563.00"
$ws.Range("H104").Value = "'70112.60"
$ws.Range("H105").Value = "'28.00"
$ws.Range("H106").Value = "'83146.71"
$ws.Range("H107").Value = "'753.00"
$ws.Range("H108").Value = "'35866.00"
$ws.Range("H109").Value = "'3660.00"
$ws.Range("H110").Value = "'78.60"
$ws.Range("H111").Value = "'10008.17"
$ws.Range("H112").Value = "'3219.00"
$ws.Range("H113").Value = "'3600.00"
$ws.Range("H114").Value = "'3252.49"
$ws.Range("H115").Value = "'5360.00"
$ws.Range("H116").Value = "'1626.80"
$ws.Range("H117").Value = "'8600.00"
$ws.Range("H118").Value = "'4248.00"
$ws.Range("H119").Value = "'10134.00"
$ws.Range("H120").Value = "'13950.00"
$ws.Range("H121").Value = "'2246.00"
$ws.Range("H122").Value = "'41019.31"
$ws.Range("H123").Value = "'5491.00"
$ws.Range("H124").Value = "'18200.00"
$ws.Range("H125").Value = "'1503.32"
$ws.Range("H126").Value = "'777.60"
$ws.Range("H127").Value = "'617.15"
$ws.Range("H128").Value = "'10000.00"
$ws.Range("H129").Value = "'7800.00"
$ws.Range("H130").Value = "'618.42"
$ws.Range("H131").Value = "'10084.12"
$ws.Range("H132").Value = "'7900.00"
$ws.Range("H133").Value = "'8000.00"
$ws.Range("H134").Value = "'5980.00"
$ws.Range("H135").Value = "'3000.00"
$ws.Range("H136").Value = "'3800.00"
$ws.Range("H137").Value = "'4532.00"
$ws.Range("H138").Value = "'11000.00"
$ws.Range("H139").Value = "'20013.00"
$ws.Range("H140").Value = "'72000.00"
$ws.Range("H141").Value = "'6000.00"
$ws.Range("H142").Value = "'7250.00"
$ws.Range("H143").Value = "'5500.00"
$ws.Range("H144").Value = "'38000.00"
$ws.Range("H145").Value = "'13240.00"
$ws.Range("H146").Value = "'189194.59"
$ws.Range("H147").Value = "'4500.00"
$ws.Range("H148").Value = "'6000.00"
$ws.Range("H149").Value = "'477.00"
$ws.Range("H150").Value = "'436.00"
$ws.Range("H151").Value = "'271.95"
$ws.Range("H152").Value = "'3393.00"
$ws.Range("H153").Value = "'728.48"
$ws.Range("H154").Value = "'31417.00"
$ws.Range("H155").Value = "'246.66"
$ws.Range("H156").Value = "'467500.00"
$ws.Range("H157").Value = "'3500.00"
$ws.Range("H158").Value = "'15000.00"
$ws.Range("H159").Value = "'20000.00"
$ws.Range("H160").Value = "'12000.00"
$ws.Range("H161").Value = "'28000.00"
$ws.Range("H162").Value = "'6000.00"
$ws.Range("H163").Value = "'48368.54"
$ws.Range("H164").Value = "'8500.00"
$ws.Range("H165").Value = "'2400.00"
$ws.Range("H166").Value = "'9998.24"
$ws.Range("H167").Value = "'8000.00"
$ws.Range("H168").Value = "'3000.00"
$ws.Range("H169").Value = "'3000.00"
$ws.Range("H170").Value = "'4500.00"
$ws.Range("H171").Value = "'7000.00"
$ws.Range("H172").Value = "'6000.00"
$ws.Range("H173").Value = "'9000.00"
$ws.Range("H174").Value = "'23424.50"
$ws.Range("H175").Value = "'2000.00"
$ws.Range("H176").Value = "'16000.00"
$ws.Range("H177").Value = "'12000.00"
$ws.Range("H178").Value = "'4500.00"
$ws.Range("H179").Value = "'12000.00"
$ws.Range("H180").Value = "'14580.00"
$ws.Range("H181").Value = "'30000.00"
$ws.Range("H182").Value = "'24300.00"
$ws.Range("H183").Value = "'18300.00"
$ws.Range("H184").Value = "'15500.00"
$ws.Range("H185").Value = "'5000.00"
$ws.Range("H186").Value = "'37880.00"
$ws.Range("H187").Value = "'4000.00"
$ws.Range("H188").Value = "'2490.00"
$ws.Range("H189").Value = "'3500.00"
$ws.Range("H190").Value = "'7950.00"
$ws.Range("H191").Value = "'54994.00"
$ws.Range("H192").Value = "'85.00"
$ws.Range("H193").Value = "'85000.00"
$ws.Range("H194").Value = "'24340.00"
$ws.Range("H195").Value = "'1996.50"
$ws.Range("H196").Value = "'1750.00"
$ws.Range("H197").Value = "'2214.00"
$ws.Range("H198").Value = "'7203.00"
$ws.Range("H199").Value = "'16990.00"
$ws.Range("H200").Value = "'4450.00"
$ws.Range("H201").Value = "'622.78"
$ws.Range("H202").Value = "'10613.22"
$ws.Range("H203").Value = "'4152.00"
$ws.Range("H204").Value = "'148200.02"
$ws.Range("H205").Value = "'23968.00"
$ws.Range("H206").Value = "'1420.00"
$ws.Range("H207").Value = "'33760.00"
$ws.Range("H208").Value = "'904.09"
$ws.Range("H209").Value = "'7670.00"
$ws.Range("H210").Value = "'1390.00"
$ws.Range("H211").Value = "'11260.00"
$ws.Range("H212").Value = "'10908.43"
$ws.Range("H213").Value = "'826.83"
$ws.Range("H214").Value = "'13350.00"
$ws.Range("H215").Value = "'10965.00"
$ws.Range("H216").Value = "'9533.00"
$ws.Range("H217").Value = "'90.00"
$ws.Range("H218").Value = "'7201.00"
$ws.Range("H219").Value = "'480.00"
$ws.Range("H220").Value = "'2320.00"
$ws.Range("H221").Value = "'900.00"
$ws.Range("H222").Value = "'1283.74"
$ws.Range("H223").Value = "'468945.60"
$ws.Range("H224").Value = "'8500.00"
$ws.Range("H225").Value = "'25000.00"
$ws.Range("H226").Value = "'25000.00"
$ws.Range("H227").Value = "'25000.00"
$ws.Range("H228").Value = "'25000.00"
$ws.Range("H229").Value = "'25000.00"
$ws.Range("H230").Value = "'50000.00"
$ws.Range("H231").Value = "'50000.00"
$ws.Range("H232").Value = "'25000.00"
$ws.Range("H233").Value = "'21950.00"
$ws.Range("H234").Value = "'7923.37"
$ws.Range("H235").Value = "'3079055.95"
$ws.Range("H236").Value = "'6200.00"
$ws.Range("H237").Value = "'216.86"
$ws.Range("H238").Value = "'2350.00"
$ws.Range("H239").Value = "'6300.00"
$ws.Range("H240").Value = "'34000.00"
$ws.Range("H241").Value = "'116190.00"
$ws.Range("H242").Value = "'122190.00"
$ws.Range("H243").Value = "'116190.00"
$ws.Range("H244").Value = "'116190.00"
$ws.Range("H245").Value = "'118670.00"
$ws.Range("H246").Value = "'116190.00"
$ws.Range("H247").Value = "'200190.00"
$ws.Range("H248").Value = "'284190.00"
$ws.Range("H249").Value = "'299910.00"
$ws.Range("H250").Value = "'116190.00"
$ws.Range("H251").Value = "'117415.00"
$ws.Range("H252").Value = "'116190.00"
$ws.Range("H253").Value = "'116190.00"
$ws.Range("H254").Value = "'116190.00"
$ws.Range("H255").Value = "'201440.00"
$ws.Range("H256").Value = "'285470.00"
$ws.Range("H257").Value = "'200190.00"
$ws.Range("H258").Value = "'116190.00"
$ws.Range("H259").Value = "'184190.00"
$ws.Range("H260").Value = "'116190.00"
$ws.Range("H261").Value = "'116190.00"
$ws.Range("H262").Value = "'122890.00"
$ws.Range("H263").Value = "'116190.00"
$ws.Range("H264").Value = "'666861.76"
$ws.Range("H265").Value = "'81102.50"
$ws.Range("H266").Value = "'60574.00"
$ws.Range("H267").Value = "'8223.00"
$ws.Range("H268").Value = "'151653.00"
$ws.Range("H269").Value = "'13800.00"
$ws.Range("H270").Value = "'5180.00"
$ws.Range("H271").Value = "'14500.00"
$ws.Range("H272").Value = "'60000.00"
$ws.Range("H273").Value = "'5000.00"
$ws.Range("H274").Value = "'2000.00"
$ws.Range("H275").Value = "'7061.00"
$ws.Range("H276").Value = "'6380.00"
$ws.Range("H277").Value = "'8800.00"
$ws.Range("H278").Value = "'2490.00"
$ws.Range("H279").Value = "'5200.00"
$ws.Range("H280").Value = "'60000.00"
$ws.Range("H281").Value = "'4300.00"
$ws.Range("H282").Value = "'16479.88"
$ws.Range("H283").Value = "'1802.04"
